# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Estado de Cuenta" for NIT 9010408211 is refreshed with newer data:
#  - The two rows belonging to the worker "YUDIS MARIA VILLERO TOVAR"
#    (document 45553717, periods 2307 and 2306) are removed from the
#    detail table, since that worker is no longer part of this report.
#  - The remaining detail rows (ROMAN MEZA and RAUL ENRIQUE VELEZ TATIS)
#    move up to take their place.
#  - The summary figures at the top of the sheet (total "Valor Mora",
#    "Cant. Trabajadores" and "Cant. Periodos") are updated to match the
#    smaller data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Remove the two detail rows for YUDIS MARIA VILLERO TOVAR (rows 16-17);
# everything below shifts up automatically.
$ws.Rows("16:17").Delete()

# Update the summary totals to reflect the reduced worker/period list.
$ws.Range("E11").Value = 39611   # VALOR MORA (total)
$ws.Range("C13").Value = 2       # Cant. Trabajadores
$ws.Range("F13").Value = 2       # Cant. Periodos
